# Use 3rd quartile instead of mean for Number_of_Inclusions (column B),
# and refresh the dependent Number_of_Inclusions_per_Nucleus (column D = B / C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => new Number_of_Inclusions value
$updates = @{
    2  = 3
    5  = 10
    6  = 16
    7  = 12
    8  = 5
    11 = 11
    14 = 3
    15 = 0
    16 = 9
    17 = 0
    18 = 14
}

foreach ($row in $updates.Keys) {
    $newB = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $newB

    $c = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value = $newB / $c
}
